$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Лист1: add two new rows (column-major fill order: A5,A6 then B5,B6) ---
$ws1.Range("A5").Value = "Корректность работы с различными входными данными"
$ws1.Range("A6").Value = "Корректность нахождения факториала числа"
$ws1.Range("B5").Value = "1; 2; 3; 4; 5"
$ws1.Range("B6").Value = "1; 3; 6; 10; 15"

# --- Лист2: drop the extra "Столбец1" column from the first table ---
$lo1 = $ws2.ListObjects.Item(1)
$lo1.ListColumns.Item(7).Delete()

# --- Лист2: populate the data for the new factorial-test table (B14:G19) ---
$ws2.Range("B14").Value = "Идентификатор теста"
$ws2.Range("C14").Value = "Аспект тестирования"
$ws2.Range("D14").Value = "Описание теста"
$ws2.Range("E14").Value = "Входные данные"
$ws2.Range("F14").Value = "Шаги выполнения"
$ws2.Range("G14").Value = "Ожидаемый результат"

$ws2.Range("B15").Value = 1
$ws2.Range("C15").Value = "Корректность работы с различными входными данными"
$ws2.Range("D15").Value = "Ввод числа"
$ws2.Range("E15").Value = 1
$ws2.Range("G15").Value = 1

$ws2.Range("B16").Value = 2
$ws2.Range("C16").Value = "Корректность работы с различными входными данными"
$ws2.Range("D16").Value = "Ввод нуля"
$ws2.Range("E16").Value = 0
$ws2.Range("G16").Value = 0

$ws2.Range("B17").Value = 3
$ws2.Range("C17").Value = "Корректность работы с различными входными данными"
$ws2.Range("D17").Value = "Ввод числа меньше нуля"
$ws2.Range("E17").Value = -1
$ws2.Range("G17").Value = "Exception"

$ws2.Range("B18").Value = 4
$ws2.Range("C18").Value = "Корректность работы с различными входными данными"
$ws2.Range("D18").Value = "Ввод большого числа"
$ws2.Range("E18").Value = 1000
$ws2.Range("G18").Value = "Большое число"

$ws2.Range("B19").Value = 5
$ws2.Range("C19").Value = "Корректность работы с различными входными данными"
$ws2.Range("D19").Value = "Ввод не числа"
$ws2.Range("E19").Value = "str"
$ws2.Range("G19").Value = "Ошибка"

# --- turn B14:G19 into a real table (ListObject), matching Таблица13 ---
$null = $ws2.ListObjects.Add(1, $ws2.Range("B14:G19"), 0, 1)
$lo2 = $ws2.ListObjects.Item(1)
$lo2.TableStyle = "TableStyleLight8"
$lo2.Name = "Таблица13"

# --- restore view state (selection / zoom) to match the saved workbook ---
$ws1.Range("B16").Select()
$ws2.Range("F29").Select()
$excel.ActiveWindow.Zoom = 130
